# Hortaliza / Zapallo weekly data update:
# A new price record (dated 2022-11-08 / serial 44873) is inserted as row 120,
# shifting all subsequent rows (old 120..199) down by one (to 121..200).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120 (shifts existing rows 120-199 down to 121-200)
$ws.Rows("120:120").Insert()

# Populate the newly inserted row with the new record's data
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 44873
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112045
$ws.Cells.Item(120, 7).Value = "Zapallo"
$ws.Cells.Item(120, 8).Value = "Camote"
$ws.Cells.Item(120, 9).Value = "1a nueva(o)"
$ws.Cells.Item(120, 10).Value = 200
$ws.Cells.Item(120, 11).Value = 1000
$ws.Cells.Item(120, 12).Value = 1100
$ws.Cells.Item(120, 13).Value = 1050
$ws.Cells.Item(120, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(120, 15).Value = "Perú"
$ws.Cells.Item(120, 16).Value = 1050
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(120, 18).Value = "Hortaliza"
